# Actualización automática del tracker
# Fills in results/profit for matches that have finished, and appends the
# newly scheduled matches for 2025-08-05.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 36: Radu Albot vs Tadeas Paroulek -> Fallo ---
$ws.Range("G36").Value = "Fallo"
$ws.Range("H36").Value = -1

# --- Row 39: Fabrizio Andaloro vs Maximus Jones -> Acierto ---
$ws.Range("G39").Value = "Acierto"
$ws.Range("H39").Value = 1.75

# --- Row 47: Alex Molcan vs Martin Krumich -> Acierto ---
$ws.Range("G47").Value = "Acierto"
$ws.Range("H47").Value = 2.4

# --- New row 49: Andrey Rublev vs Taylor Fritz ---
$ws.Range("A49").Value = 14265597
$b49 = $ws.Cells.Item(49, 2)
$b49.Value = "'2025-08-05"
$b49.Style = "Normal"
$ws.Range("C49").Value = "Andrey Rublev"
$ws.Range("D49").Value = "Taylor Fritz"
$ws.Range("E49").Value = "Gana Andrey Rublev"
$ws.Range("F49").Value = 2.75

# --- New row 50: Madison Keys vs Clara Tauson ---
$ws.Range("A50").Value = 14266319
$b50 = $ws.Cells.Item(50, 2)
$b50.Value = "'2025-08-05"
$b50.Style = "Normal"
$ws.Range("C50").Value = "Madison Keys"
$ws.Range("D50").Value = "Clara Tauson"
$ws.Range("E50").Value = "Gana Clara Tauson"
$ws.Range("F50").Value = 2.3

$wb.Save()
